$wb = $excel.ActiveWorkbook

# "table attribute" sheet (sheet1) gains a new column B header: "table description"
$ws1 = $wb.Worksheets.Item("table attribute")
$ws1.Range("B1").Value = "table description"

# copy A1's formatting (style index) onto B1 so it matches the header style
$ws1.Range("A1").Copy()
$ws1.Range("B1").PasteSpecial(-4122)
$ws1.Range("B1").Value = "table description"

# The user switched focus to the "table attribute" tab and selected B4 there
$ws1.Activate()
$ws1.Range("B4").Select()
